$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Update cell A11 text
$ws.Range("A11").Value = "Function Information"

# 3. Update cell B21 value
$ws.Range("B21").Value = 0.2499825759175085

# 4. Update number format code for the data column (numFmtId 166: "0.000" -> "###0.000")
$ws.Range("B27:B36").NumberFormat = "###0.000"
